$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")

# New order rows appended to the Orders sheet (rows 9-20): a pending
# order for "Sakamoto Days" volumes (rows 9-19, user "assa") and a
# completed order for the same manga (row 20, user "vxbvnb").
$rows = @(
    @("ياسين المريني", "assa",   "43215665454765", "PENDING", "2025-07-31T16:43:59.649Z", "Sakamoto Days", 10, 1, 118.8, 118.8,  1401.3),
    @("ياسين المريني", "assa",   "43215665454765", "PENDING", "2025-07-31T16:43:59.649Z", "Sakamoto Days", 6,  1, 118.8, 118.8,  1401.3),
    @("ياسين المريني", "assa",   "43215665454765", "PENDING", "2025-07-31T16:43:59.649Z", "Sakamoto Days", 5,  1, 118.8, 118.8,  1401.3),
    @("ياسين المريني", "assa",   "43215665454765", "PENDING", "2025-07-31T16:43:59.649Z", "Sakamoto Days", 4,  1, 118.8, 118.8,  1401.3),
    @("ياسين المريني", "assa",   "43215665454765", "PENDING", "2025-07-31T16:43:59.649Z", "Sakamoto Days", 3,  1, 118.8, 118.8,  1401.3),
    @("ياسين المريني", "assa",   "43215665454765", "PENDING", "2025-07-31T16:43:59.649Z", "Sakamoto Days", 2,  1, 118.8, 118.8,  1401.3),
    @("ياسين المريني", "assa",   "43215665454765", "PENDING", "2025-07-31T16:43:59.649Z", "Sakamoto Days", 1,  1, 118.8, 118.8,  1401.3),
    @("ياسين المريني", "assa",   "43215665454765", "PENDING", "2025-07-31T16:43:59.649Z", "Sakamoto Days", 9,  1, 118.8, 118.8,  1401.3),
    @("ياسين المريني", "assa",   "43215665454765", "PENDING", "2025-07-31T16:43:59.649Z", "Sakamoto Days", 8,  2, 118.8, 237.6,  1401.3),
    @("ياسين المريني", "assa",   "43215665454765", "PENDING", "2025-07-31T16:43:59.649Z", "Sakamoto Days", 7,  1, 118.8, 118.8,  1401.3),
    @("ياسين المريني", "assa",   "43215665454765", "PENDING", "2025-07-31T16:43:59.649Z", "Jujutsu Kaisen", 5,  1, 94.5,  94.5,   1401.3),
    @("ياسين المريني", "vxbvnb", "435436465565",   "PENDING", "2025-07-31T16:44:54.370Z", "Sakamoto Days", 8,  8, 118.8, 950.4,  950.4)
)

$startRow = 9
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    for ($c = 0; $c -lt $data.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c + 1)
        $val = $data[$c]
        if ($c -eq 2) {
            # Phone numbers are digit strings; writing them through .Value
            # makes Excel coerce them to a number, so instead enter them as
            # a text formula and flatten it to a plain cached text value via
            # copy / paste-special-values (keeps cell style untouched).
            $cell.Formula = '="' + $val + '"'
            $cell.Copy()
            $cell.PasteSpecial(-4163)
        } else {
            $cell.Value = $val
        }
    }
}

$excel.CutCopyMode = 0
